$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AF2").Value = 9.5
$ws.Range("AI2").Value = 451
$ws.Range("I2").Value = 11
$ws.Range("L2").Value = 8.5
$ws.Range("N2").Value = 12
$ws.Range("Z2").Value = 6
$ws.Range("AB3").Value = 17
$ws.Range("AK3").Value = 19
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("AP5").Value = 1.98
$ws.Range("AQ5").Value = 1.88
$ws.Range("AR5").Value = 4.2
$ws.Range("AS5").Value = 1.23
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.48
$ws.Range("S5").Value = 5.5
$ws.Range("T5").Value = 1.14
$ws.Range("AR6").Value = 4.8
$ws.Range("AJ8").Value = 9
$ws.Range("J8").Value = 2.6
$ws.Range("M8").Value = 1.14
$ws.Range("N8").Value = 5.5
$ws.Range("AI9").Value = 251
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
$ws.Range("Q9").Value = 2.1
$ws.Range("R9").Value = 1.7
$ws.Range("Y9").Value = 8
$ws.Range("Z9").Value = 12
$ws.Range("O11").Value = 1.3
$ws.Range("P11").Value = 3.4
$ws.Range("AA13").Value = 11
$ws.Range("AB13").Value = 32
$ws.Range("AC13").Value = 26
$ws.Range("AD13").Value = 40
$ws.Range("AF13").Value = 6.7
$ws.Range("AJ13").Value = 7.6
$ws.Range("AK13").Value = 12.5
$ws.Range("AL13").Value = 10.25
$ws.Range("AM13").Value = 28
$ws.Range("AN13").Value = 23
$ws.Range("G13").Value = 2.72
$ws.Range("H13").Value = 3.35
$ws.Range("I13").Value = 2.47
$ws.Range("J13").Value = 3.25
$ws.Range("L13").Value = 3
$ws.Range("Y13").Value = 7.9
$ws.Range("Z13").Value = 14
$ws.Range("AI17").Value = 351
$ws.Range("AK17").Value = 8.5
$ws.Range("AM17").Value = 15
$ws.Range("AN17").Value = 15
$ws.Range("G17").Value = 3.8
$ws.Range("H17").Value = 3.7
$ws.Range("I17").Value = 1.9
$ws.Range("Q17").Value = 2
$ws.Range("AB18").Value = 12.5
$ws.Range("AC18").Value = 13
$ws.Range("AF18").Value = 7.1
$ws.Range("AG18").Value = 16
$ws.Range("AH18").Value = 75
$ws.Range("AJ18").Value = 13
$ws.Range("AK18").Value = 28
$ws.Range("AL18").Value = 15.5
$ws.Range("AM18").Value = 90
$ws.Range("AN18").Value = 45
$ws.Range("AO18").Value = 50
$ws.Range("G18").Value = 1.65
$ws.Range("H18").Value = 3.65
$ws.Range("I18").Value = 4.75
$ws.Range("J18").Value = 2.22
$ws.Range("L18").Value = 4.9
$ws.Range("P18").Value = 3.15
$ws.Range("Q18").Value = 1.78
$ws.Range("R18").Value = 1.82
$ws.Range("S18").Value = 2.82
$ws.Range("T18").Value = 1.33
$ws.Range("W18").Value = 1.78
$ws.Range("X18").Value = 1.83
$ws.Range("Y18").Value = 6.9
$ws.Range("Z18").Value = 7.7
$ws.Range("AB19").Value = 26
$ws.Range("AC19").Value = 18
$ws.Range("AD19").Value = 23
$ws.Range("AE19").Value = 13.5
$ws.Range("AF19").Value = 7.2
$ws.Range("AG19").Value = 12.5
$ws.Range("AH19").Value = 45
$ws.Range("AI19").Value = 250
$ws.Range("AJ19").Value = 10.75
$ws.Range("AK19").Value = 14.5
$ws.Range("AN19").Value = 19
$ws.Range("AO19").Value = 24
$ws.Range("H19").Value = 3.6
$ws.Range("I19").Value = 2.52
$ws.Range("J19").Value = 2.92
$ws.Range("K19").Value = 2.25
$ws.Range("L19").Value = 3
$ws.Range("O19").Value = 1.2
$ws.Range("P19").Value = 3.6
$ws.Range("Q19").Value = 1.6
$ws.Range("R19").Value = 2.05
$ws.Range("S19").Value = 2.45
$ws.Range("T19").Value = 1.42
$ws.Range("W19").Value = 1.53
$ws.Range("X19").Value = 2.18
$ws.Range("Y19").Value = 10.5
$ws.Range("Z19").Value = 13.5
$ws.Range("N20").Value = 9
$ws.Range("O20").Value = 1.33
$ws.Range("P20").Value = 3.25
$ws.Range("AB22").Value = 26
$ws.Range("AC22").Value = 19
$ws.Range("AJ22").Value = 11
$ws.Range("AL22").Value = 10
$ws.Range("AN22").Value = 19
$ws.Range("G22").Value = 2.63
$ws.Range("I22").Value = 2.45
$ws.Range("J22").Value = 3.2
$ws.Range("L22").Value = 3.1
$ws.Range("AC24").Value = 21
$ws.Range("J24").Value = 3.75
$ws.Range("L24").Value = 2.3
$ws.Range("O24").Value = 1.1
$ws.Range("P24").Value = 7
$ws.Range("S24").Value = 1.83
$ws.Range("T24").Value = 1.83
$ws.Range("AJ28").Value = 8
$ws.Range("G28").Value = 2.38
$ws.Range("I28").Value = 2.75
$ws.Range("L28").Value = 3.6
$ws.Range("M28").Value = 1.07
$ws.Range("N28").Value = 9
$ws.Range("O28").Value = 1.36
$ws.Range("P28").Value = 3
$ws.Range("Z28").Value = 11
$ws.Range("M29").Value = 1.06
$ws.Range("N29").Value = 10
$ws.Range("Q29").Value = 2.1
$ws.Range("R29").Value = 1.7
$ws.Range("S29").Value = 3.75
$ws.Range("T29").Value = 1.25
$ws.Range("G31").Value = 1.7
$ws.Range("H31").Value = 3.5
$ws.Range("I31").Value = 5
$ws.Range("M31").Value = 1.05
$ws.Range("N31").Value = 8.5
$ws.Range("Q31").Value = 1.88
$ws.Range("R31").Value = 1.93
$ws.Range("AA32").Value = 10.75
$ws.Range("AC32").Value = 25
$ws.Range("AD32").Value = 30
$ws.Range("AJ32").Value = 8.75
$ws.Range("AN32").Value = 16.5
$ws.Range("AO32").Value = 24
$ws.Range("G32").Value = 3
$ws.Range("I32").Value = 2.15
$ws.Range("J32").Value = 3.6
$ws.Range("L32").Value = 2.75
$ws.Range("P32").Value = 3.5
$ws.Range("S32").Value = 2.82
$ws.Range("T32").Value = 1.38
$ws.Range("W32").Value = 1.62
$ws.Range("X32").Value = 2.15
$ws.Range("Y32").Value = 10.25
$ws.Range("AA34").Value = 12.5
$ws.Range("AB34").Value = 5.7
$ws.Range("AC34").Value = 12
$ws.Range("AD34").Value = 45
$ws.Range("AE34").Value = 17.5
$ws.Range("AF34").Value = 18.5
$ws.Range("AH34").Value = 250
$ws.Range("AJ34").Value = 80
$ws.Range("AN34").Value = 700
$ws.Range("AO34").Value = 350
$ws.Range("G34").Value = 1.09
$ws.Range("H34").Value = 7.3
$ws.Range("I34").Value = 25
$ws.Range("J34").Value = 1.39
$ws.Range("K34").Value = 3
$ws.Range("L34").Value = 15.5
$ws.Range("O34").Value = 1.12
$ws.Range("P34").Value = 5.8
$ws.Range("Q34").Value = 1.37
$ws.Range("R34").Value = 2.62
$ws.Range("S34").Value = 1.9
$ws.Range("T34").Value = 1.72
$ws.Range("U34").Value = 1.19
$ws.Range("V34").Value = 4.25
$ws.Range("W34").Value = 2.47
$ws.Range("X34").Value = 1.42
$ws.Range("Y34").Value = 8.5
$ws.Range("Z34").Value = 5.8
